# Append new conversation-log rows to the "대화기록" (chat log) sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("대화기록")

$rows = @(
    @("2026-02-26", "11:11:40", "Q: 유리컵 냉장고에 있던데? | A: 유리컵을 냉장고에 보관하는 것은 일반적으로 음료를 차갑게 유지하기 위해서입니다. 특히 여름철에 시원한 음료를 마시고 싶을 때 유용할 수 있습니다. 하지만 유리컵은 온도 변화에 따라 깨질 수 있으니 주의해야 합니다."),
    @("2026-02-26", "11:11:50", "Q: 유리컵 어디있어? | A: 컵은 선반에 있어요."),
    @("2026-02-26", "13:54:35", "Q: 유리컵 어디있어? | A: 컵은 선반에 있어요."),
    @("2026-02-26", "13:56:08", "Q: 유리컵 어디있어? | A: 컵은 선반에 있어요."),
    @("2026-02-26", "14:01:08", "Q: 유리컵 어디있어? | A: 컵은 선반에 있어요."),
    @("2026-02-26", "14:07:16", "Q: 유리컵 어디있어? | A: 컵은 선반에 있어요."),
    @("2026-02-26", "14:08:55", "Q: 유리컵 어디있어? | A: 컵은 선반에 있어요.")
)

$startRow = 14

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i

    # The "date" column holds plain text like "2026-02-26" (same as the
    # rows above it), but Excel auto-converts that text to a real date
    # value on assignment. Force it to stay text, then re-apply the
    # (unstyled) look of the existing data rows so no stray number
    # format / style index is left behind on the cell.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 1).Style = $ws.Cells.Item(2, 1).Style

    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
